# Update beverage description text on the active worksheet to match the
# revised copy in the menu template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14").Value = "Strong & popular"
$ws.Range("B15").Value = "Strong & popular"
$ws.Range("B16").Value = "Strong & popular"
$ws.Range("B17").Value = "Easy drinking"
$ws.Range("B18").Value = "Easy drinking"
$ws.Range("B19").Value = "Easy drinking"
$ws.Range("B20").Value = "Easy drinking"
$ws.Range("B27").Value = "A powerful blend of five white spirits, citrus, and cola. Our strongest classic."
$ws.Range("B60").Value = "Shareable, fun & flavorful wine cocktails"
